# Fix[ITE2]: general fixes on table creations
#
# Applies the cell-content corrections from the "modelo-relacional" sheet:
# several attribute-name cells that used to repeat a table's own entity
# name (e.g. "pedido", "sucursal", "producto", "Compra", "proveedor") are
# renamed to their PK-style id form ("idPedido", "idSucursal",
# "idProducto", "idCompra", "idProveedor"); two overlong / stale check
# constraint comments are trimmed; a missing "O" (optional) marker cell is
# added; and the view is left scrolled/zoomed/selected near the edited
# area.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Table "Producto" (rows 8-9): trim stale/overlong check constraint ---
$ws.Range("D9").Value = "NN, CK in F"

# --- Table "Review_Pedido" (rows 8-9) ---
$ws.Range("G8").Value = "idProducto"

# --- Table "Productos_Compra" (rows 17-19) ---
$ws.Range("L18").Value = "idProveedor"

# --- Table "Pedidos_Proveedor" (rows 3-4 header block) ---
$ws.Range("H3").Value = "idPedido"

# --- Table "Medición_Producto" (row 14): stale check-constraint list ---
$ws.Range("L14").Value = "NN, CK in  "

$ws.Range("I3").Value = "idSucursal"

$ws.Range("H18").Value = "idCompra"

# New "O" (optional) marker cell next to the Productos_Compra block.
$ws.Range("J18").Value = "O"

$ws.Range("I18").Value = "idProducto"
$ws.Range("K18").Value = "idPedido"

# Row 14 visibly shrinks now that its longest cell holds much shorter text.
$ws.Rows.Item(14).RowHeight = 17

# --- View state: scrolled/zoomed near the newly-edited J18 cell ---
$ws.Range("J18").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 75
